$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain literal text assignments (dates, hPa, km/h, degC, MJ/m2, etc.)
# These do not get misinterpreted as numbers/dates by Excel's smart-entry parser.
$ws.Range("E2").Value = "2026-02-25 05:48:04"
$ws.Range("E3").Value = "2026-02-25 05:48:06"
$ws.Range("E4").Value = "2026-02-25 05:48:08"
$ws.Range("J4").Value = "1019.7 hPa"
$ws.Range("N4").Value = "1.2 °C 5:13 TU"
$ws.Range("O4").Value = "2.6 °C"
$ws.Range("E5").Value = "2026-02-25 05:48:11"
$ws.Range("E6").Value = "2026-02-25 05:48:13"
$ws.Range("J6").Value = "1019.3 hPa"
$ws.Range("O6").Value = "8.4 °C"
$ws.Range("E7").Value = "2026-02-25 05:48:15"
$ws.Range("L7").Value = "13.0 km/h - 75º 5:05 TU"
$ws.Range("N7").Value = "10.0 °C 5:22 TU"
$ws.Range("E8").Value = "2026-02-25 05:48:18"
$ws.Range("J8").Value = "1018.1 hPa"
$ws.Range("K8").Value = "-0.1 MJ/m2"
$ws.Range("N8").Value = "12.1 °C 5:22 TU"
$ws.Range("O8").Value = "14.8 °C"
$ws.Range("E9").Value = "2026-02-25 05:48:20"
$ws.Range("N9").Value = "3.2 °C 5:25 TU"
$ws.Range("O9").Value = "5.0 °C"
$ws.Range("E10").Value = "2026-02-25 05:48:22"
$ws.Range("L10").Value = "5.8 km/h - 99º 5:04 TU"
$ws.Range("M10").Value = "6.5 °C 5:11 TU"
$ws.Range("O10").Value = "4.2 °C"
$ws.Range("E11").Value = "2026-02-25 05:48:25"
$ws.Range("O11").Value = "2.9 °C"
$ws.Range("E12").Value = "2026-02-25 05:48:27"
$ws.Range("O12").Value = "5.4 °C"
$ws.Range("E13").Value = "2026-02-25 05:48:29"
$ws.Range("J13").Value = "1026.8 hPa"
$ws.Range("N13").Value = "-3.8 °C 5:22 TU"
$ws.Range("O13").Value = "-1.7 °C"
$ws.Range("E14").Value = "2026-02-25 05:48:31"
$ws.Range("N14").Value = "2.8 °C 5:29 TU"
$ws.Range("O14").Value = "5.3 °C"
$ws.Range("E15").Value = "2026-02-25 05:48:34"
$ws.Range("N15").Value = "3.9 °C 5:29 TU"
$ws.Range("O15").Value = "5.4 °C"
$ws.Range("E16").Value = "2026-02-25 05:48:36"
$ws.Range("O16").Value = "3.6 °C"
$ws.Range("E17").Value = "2026-02-25 05:48:38"
$ws.Range("O17").Value = "8.9 °C"
$ws.Range("E18").Value = "2026-02-25 05:48:40"
$ws.Range("J18").Value = "1019.6 hPa"
$ws.Range("N18").Value = "4.7 °C 5:21 TU"
$ws.Range("O18").Value = "6.1 °C"
$ws.Range("E19").Value = "2026-02-25 05:48:42"
$ws.Range("O19").Value = "9.7 °C"
$ws.Range("E20").Value = "2026-02-25 05:48:45"
$ws.Range("N20").Value = "2.0 °C 5:25 TU"
$ws.Range("E21").Value = "2026-02-25 05:48:47"
$ws.Range("J21").Value = "1023.6 hPa"
$ws.Range("N21").Value = "1.4 °C 5:10 TU"
$ws.Range("O21").Value = "3.5 °C"
$ws.Range("E22").Value = "2026-02-25 05:48:49"
$ws.Range("O22").Value = "1.4 °C"
$ws.Range("E23").Value = "2026-02-25 05:48:52"
$ws.Range("M23").Value = "4.2 °C 5:26 TU"
$ws.Range("O23").Value = "3.4 °C"
$ws.Range("E24").Value = "2026-02-25 05:48:54"
$ws.Range("L24").Value = "8.3 km/h - 39º 5:17 TU"
$ws.Range("E25").Value = "2026-02-25 05:48:56"
$ws.Range("N25").Value = "2.2 °C 5:26 TU"
$ws.Range("O25").Value = "3.1 °C"
$ws.Range("E26").Value = "2026-02-25 05:48:59"
$ws.Range("N26").Value = "7.8 °C 5:11 TU"
$ws.Range("O26").Value = "9.1 °C"
$ws.Range("E27").Value = "2026-02-25 05:49:01"
$ws.Range("E28").Value = "2026-02-25 05:49:03"
$ws.Range("J28").Value = "1020.6 hPa"
$ws.Range("O28").Value = "3.6 °C"
$ws.Range("E29").Value = "2026-02-25 05:49:05"
$ws.Range("M29").Value = "11.6 °C 5:29 TU"
$ws.Range("O29").Value = "9.3 °C"
$ws.Range("E30").Value = "2026-02-25 05:49:07"
$ws.Range("J30").Value = "1019.5 hPa"
$ws.Range("N30").Value = "6.2 °C 5:29 TU"
$ws.Range("O30").Value = "7.5 °C"
$ws.Range("E31").Value = "2026-02-25 05:49:10"
$ws.Range("J31").Value = "1018.7 hPa"
$ws.Range("E32").Value = "2026-02-25 05:49:12"
$ws.Range("N32").Value = "0.4 °C 5:20 TU"
$ws.Range("O32").Value = "2.2 °C"
$ws.Range("E33").Value = "2026-02-25 05:49:14"
$ws.Range("J33").Value = "1024.1 hPa"
$ws.Range("N33").Value = "0.6 °C 5:29 TU"
$ws.Range("O33").Value = "2.2 °C"
$ws.Range("E34").Value = "2026-02-25 05:49:16"
$ws.Range("E35").Value = "2026-02-25 05:49:19"
$ws.Range("N35").Value = "8.2 °C 5:19 TU"
$ws.Range("O35").Value = "9.8 °C"
$ws.Range("E36").Value = "2026-02-25 05:49:21"
$ws.Range("J36").Value = "1019.2 hPa"
$ws.Range("E37").Value = "2026-02-25 05:49:23"
$ws.Range("J37").Value = "1024.2 hPa"
$ws.Range("O37").Value = "0.9 °C"
$ws.Range("E38").Value = "2026-02-25 05:49:25"
$ws.Range("N38").Value = "3.2 °C 5:29 TU"
$ws.Range("E39").Value = "2026-02-25 05:49:28"
$ws.Range("E40").Value = "2026-02-25 05:49:30"
$ws.Range("N40").Value = "0.0 °C 5:17 TU"
$ws.Range("O40").Value = "1.4 °C"
$ws.Range("E41").Value = "2026-02-25 05:49:32"
$ws.Range("J41").Value = "1018.6 hPa"
$ws.Range("E42").Value = "2026-02-25 05:49:34"
$ws.Range("O42").Value = "8.4 °C"
$ws.Range("E43").Value = "2026-02-25 05:49:37"
$ws.Range("E44").Value = "2026-02-25 05:49:39"
$ws.Range("O44").Value = "-0.4 °C"
$ws.Range("E45").Value = "2026-02-25 05:49:41"
$ws.Range("N45").Value = "4.8 °C 5:05 TU"
$ws.Range("E46").Value = "2026-02-25 05:49:44"
$ws.Range("J46").Value = "1019.5 hPa"
$ws.Range("N46").Value = "2.0 °C 5:29 TU"
$ws.Range("O46").Value = "3.4 °C"

# Percentage-looking text (e.g. "96%") would be auto-converted to a numeric
# percentage by a plain .Value assignment (Excel's AutoPercentEntry behaviour),
# which also reassigns the cell a new number-format style. To keep these cells
# as literal text (same style, same shown characters) we build the text via a
# formula (a quoted string literal always evaluates to Text) and then collapse
# the formula down to a static value with a values-only paste.
$c = $ws.Range("H6")
$c.Formula = '="96%"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("H8")
$c.Formula = '="49%"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("H17")
$c.Formula = '="26%"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("H21")
$c.Formula = '="76%"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("H25")
$c.Formula = '="26%"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("H26")
$c.Formula = '="39%"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("H33")
$c.Formula = '="67%"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("H34")
$c.Formula = '="53%"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("H35")
$c.Formula = '="36%"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("H45")
$c.Formula = '="47%"'
$c.Copy()
$c.PasteSpecial(-4163)
$excel.CutCopyMode = $false
